# Fix labels and restructure from_data_to_dashboard.
#
# - Moves the "base_mental_health" row up from row 9 to row 7 (age_bins and
#   sector shift down to rows 8 and 9 respectively).
# - Fixes the nice_name_english label for base_mental_health from the
#   duplicated "Mental Health Before" to "Mental Health".
# - Applies a blue font colour to the (now relocated) age_bins and sector
#   rows (8-9) to flag them as restructured.
# - Widens column A slightly to fit the longer labels.
# - Updates the active selection to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: base_mental_health (moved up from old row 9; nice_name fixed)
$ws.Range("A7").Value = "base_mental_health"
$ws.Range("B7").Value = "Mental Health Before"
$ws.Range("C7").Value = "Mental Health"
$ws.Range("D7").Value = "Background Variables"
$ws.Range("E7").Value = "Background Overview"

# Row 8: age_bins (moved down from old row 7)
$ws.Range("A8").Value = "age_bins"
$ws.Range("B8").Value = "Fine Grained Age Group"
$ws.Range("C8").Value = "Age"
$ws.Range("D8").Value = "Background Variables"
$ws.Range("E8").Value = "Background Correlation"

# Row 9: sector (moved down from old row 8)
$ws.Range("A9").Value = "sector"
$ws.Range("B9").Value = "Sector"
$ws.Range("C9").Value = "Sector"
$ws.Range("D9").Value = "Background Variables"
$ws.Range("E9").Value = "Background Correlation"

# Highlight the restructured rows (age_bins / sector) in blue.
$ws.Range("A8:E9").Font.Color = 10108961

# Column A needs to be a bit wider to fit "base_mental_health" etc.
$ws.Columns.Item(1).ColumnWidth = 23

# Update the active selection/cursor position.
$ws.Range("C8").Select()
